# Generate Report for Handback
#
# 1. The generic "Ready for handoff" status (Overview!E8/F8, zh-cn!C8, de-de!C8)
#    is replaced everywhere with "Handback transform failed" because the
#    handback transform for file 97aa4c9e-7306-4c55-a2d0-ca391f4ca931 failed.
# 2. The per-locale "Error Detail" column (R) on the zh-cn and de-de sheets is
#    populated with the reason: the handback file name returned by the
#    translator did not match the expected handoff file name.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E8").Value = "Handback transform failed"
$wsOverview.Range("F8").Value = "Handback transform failed"
$wsZhCn.Range("C8").Value = "Handback transform failed"
$wsDeDe.Range("C8").Value = "Handback transform failed"

$wsZhCn.Range("R8").Value = "Handback file name: vrejjin3.3g2 is different with handoff file name: 97aa4c9e-7306-4c55-a2d0-ca391f4ca931.c680500a63c76147ed4b0fbe248e08787a686f41.zh-cn."
$wsDeDe.Range("R8").Value = "Handback file name: vrejjin3.3g2 is different with handoff file name: 97aa4c9e-7306-4c55-a2d0-ca391f4ca931.c680500a63c76147ed4b0fbe248e08787a686f41.de-de."
